$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 11209
$ws.Range("J32").Value = 15224.875
$ws.Range("L32").Value = 15224.875
$ws.Range("N32").Value = -15876.875
$ws.Range("H80").Value = 910.58826
$ws.Range("I80").Value = 1366.6666
$ws.Range("J80").Value = 661.8182
$ws.Range("K80").Value = 4099.9998
$ws.Range("L80").Value = 1985.4546
$ws.Range("M80").Value = -3101.9998
$ws.Range("N80").Value = -3981.4546
$ws.Range("H83").Value = 910.58826
$ws.Range("I83").Value = 1366.6666
$ws.Range("J83").Value = 661.8182
$ws.Range("K83").Value = 12299.9994
$ws.Range("L83").Value = 5956.3638
$ws.Range("M83").Value = -7307.999400000001
$ws.Range("N83").Value = -15940.3638
$ws.Range("H112").Value = 2116.6316
$ws.Range("J112").Value = 2236.2354
$ws.Range("L112").Value = 6708.706200000001
$ws.Range("N112").Value = -8924.706200000001
$ws.Range("H129").Value = 876.3043
$ws.Range("J129").Value = 941.5333000000001
$ws.Range("L129").Value = 2824.5999
$ws.Range("N129").Value = -12824.5999
$ws.Range("H137").Value = 3091.4348
$ws.Range("I137").Value = 2898.6428
$ws.Range("J137").Value = 3391.3333
$ws.Range("K137").Value = 8695.928400000001
$ws.Range("L137").Value = 10173.9999
$ws.Range("M137").Value = -6145.928400000001
$ws.Range("N137").Value = -15273.9999
$ws.Range("H138").Value = 2745.0227
$ws.Range("I138").Value = 2508.75
$ws.Range("J138").Value = 2768.65
$ws.Range("K138").Value = 7526.25
$ws.Range("L138").Value = 8305.950000000001
$ws.Range("M138").Value = -2386.25
$ws.Range("N138").Value = -18585.95

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1891.6522
$ws.Range("I2").Value = 1628.625
$ws.Range("J2").Value = 2492.8572
$ws.Range("K2").Value = 1628.625
$ws.Range("L2").Value = 2492.8572
$ws.Range("M2").Value = -1515.625
$ws.Range("N2").Value = -2718.8572
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 0
$ws.Range("K4").Value = 0
$ws.Range("M4").ClearContents()
$ws.Range("H32").Value = 11961.69
$ws.Range("I32").Value = 9014.963
$ws.Range("K32").Value = 9014.963
$ws.Range("M32").Value = -8727.963
$ws.Range("H45").Value = 1402.375
$ws.Range("I45").Value = 1246.2858
$ws.Range("J45").Value = 1523.7778
$ws.Range("K45").Value = 1246.2858
$ws.Range("L45").Value = 1523.7778
$ws.Range("M45").Value = -869.2858000000001
$ws.Range("N45").Value = -2277.7778
$ws.Range("H116").Value = 1891.6522
$ws.Range("I116").Value = 1628.625
$ws.Range("J116").Value = 2492.8572
$ws.Range("K116").Value = 1628.625
$ws.Range("L116").Value = 2492.8572
$ws.Range("M116").Value = 665.375
$ws.Range("N116").Value = -7080.8572

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1891.6522
$ws.Range("I3").Value = 1628.625
$ws.Range("J3").Value = 2492.8572
$ws.Range("K3").Value = 1628.625
$ws.Range("L3").Value = 2492.8572
$ws.Range("M3").Value = -1514.625
$ws.Range("N3").Value = -2720.8572

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 229.0625
$ws.Range("I7").Value = 124.63636
$ws.Range("J7").Value = 458.8
$ws.Range("K7").Value = 124.63636
$ws.Range("L7").Value = 458.8
$ws.Range("M7").Value = -11.63636
$ws.Range("N7").Value = -684.8
$ws.Range("H31").Value = 2253.0286
$ws.Range("I31").Value = 2186.625
$ws.Range("J31").Value = 2961.3333
$ws.Range("K31").Value = 2186.625
$ws.Range("L31").Value = 2961.3333
$ws.Range("M31").Value = -1891.625
$ws.Range("N31").Value = -3551.3333
$ws.Range("H34").Value = 2253.0286
$ws.Range("I34").Value = 2186.625
$ws.Range("J34").Value = 2961.3333
$ws.Range("K34").Value = 2186.625
$ws.Range("L34").Value = 2961.3333
$ws.Range("M34").Value = -1984.625
$ws.Range("N34").Value = -3365.3333
$ws.Range("H62").Value = 9526111
$ws.Range("I62").Value = 2396.842
$ws.Range("K62").Value = 2396.842
$ws.Range("M62").Value = -1772.842
$ws.Range("H65").Value = 9526111
$ws.Range("I65").Value = 2396.842
$ws.Range("K65").Value = 11984.21
$ws.Range("M65").Value = -8864.210000000001
$ws.Range("H134").Value = 10205437
$ws.Range("I134").Value = 1391.5625
$ws.Range("K134").Value = 4174.6875
$ws.Range("M134").Value = -1639.6875

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H15").Value = 850
$ws.Range("J15").Value = 300
$ws.Range("L15").Value = 900
$ws.Range("N15").Value = -1180
$ws.Range("H113").Value = 690.1212
$ws.Range("I113").Value = 579.5
$ws.Range("J113").Value = 738.2174
$ws.Range("K113").Value = 1738.5
$ws.Range("L113").Value = 2214.6522
$ws.Range("M113").Value = 431.5
$ws.Range("N113").Value = -6554.6522
$ws.Range("H131").Value = 38520300
$ws.Range("J131").Value = 76242.75
$ws.Range("L131").Value = 228728.25
$ws.Range("N131").Value = -238808.25
$ws.Range("H138").Value = 3123.4348
$ws.Range("I138").Value = 3543.9
$ws.Range("J138").Value = 2800
$ws.Range("K138").Value = 10631.7
$ws.Range("L138").Value = 8400
$ws.Range("M138").Value = -5491.700000000001
$ws.Range("N138").Value = -18680

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H62").Value = 15085
$ws.Range("J62").Value = 15085
$ws.Range("L62").Value = 15085
$ws.Range("N62").Value = -16457
$ws.Range("H65").Value = 15085
$ws.Range("J65").Value = 15085
$ws.Range("L65").Value = 45255
$ws.Range("N65").Value = -52119
$ws.Range("H80").Value = 3580.8333
$ws.Range("I80").Value = 1798
$ws.Range("K80").Value = 1798
$ws.Range("M80").Value = -800
$ws.Range("H83").Value = 3580.8333
$ws.Range("I83").Value = 1798
$ws.Range("K83").Value = 8990
$ws.Range("M83").Value = -3998

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1150.6296
$ws.Range("I7").Value = 926.6818
$ws.Range("J7").Value = 2136
$ws.Range("K7").Value = 926.6818
$ws.Range("L7").Value = 2136
$ws.Range("M7").Value = -814.6818
$ws.Range("N7").Value = -2360
$ws.Range("H40").Value = 2175.3684
$ws.Range("I40").Value = 2068.7334
$ws.Range("K40").Value = 2068.7334
$ws.Range("M40").Value = -1932.7334
$ws.Range("H55").Value = 812.75
$ws.Range("I55").Value = 417
$ws.Range("J55").Value = 2000
$ws.Range("K55").Value = 417
$ws.Range("L55").Value = 2000
$ws.Range("M55").Value = -244
$ws.Range("N55").Value = -2346
$ws.Range("H122").Value = 20836082
$ws.Range("I122").Value = 25002598
$ws.Range("K122").Value = 75007794
$ws.Range("M122").Value = -75005344
$ws.Range("H126").Value = 1150.6296
$ws.Range("I126").Value = 926.6818
$ws.Range("J126").Value = 2136
$ws.Range("K126").Value = 2780.0454
$ws.Range("L126").Value = 6408
$ws.Range("M126").Value = -310.0454
$ws.Range("N126").Value = -11348

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H69").Value = 25000
$ws.Range("J69").Value = 25000
$ws.Range("L69").Value = 25000
$ws.Range("N69").Value = -26498
$ws.Range("H72").Value = 25000
$ws.Range("J72").Value = 25000
$ws.Range("L72").Value = 75000
$ws.Range("N72").Value = -82488
$ws.Range("H96").Value = 0
$ws.Range("I96").Value = 0
$ws.Range("J96").Value = 0
$ws.Range("K96").Value = 0
$ws.Range("L96").Value = 0
$ws.Range("M96").ClearContents()
$ws.Range("N96").ClearContents()
$ws.Range("H122").Value = 9617679
$ws.Range("I122").Value = 12502323
$ws.Range("K122").Value = 37506969
$ws.Range("M122").Value = -37504519
$ws.Range("H132").Value = 3141.5715
$ws.Range("I132").Value = 3129.524
$ws.Range("J132").Value = 3213.8572
$ws.Range("K132").Value = 9388.572
$ws.Range("L132").Value = 9641.571599999999
$ws.Range("M132").Value = -6858.572
